# Apply updated CSV/XLSX values to the JFK, Regular, and Others sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: JFK ---
$ws = $wb.Worksheets.Item("JFK")
$ws.Range("C2").Value = 38760
$ws.Range("D2").Value = 709282.13
$ws.Range("C3").Value = 21115
$ws.Range("D3").Value = 387399.07
$ws.Range("C4").Value = 46591
$ws.Range("D4").Value = 849272.37
$ws.Range("C5").Value = 21762
$ws.Range("D5").Value = 396681.62
$ws.Range("C6").Value = 77375
$ws.Range("D6").Value = 1408245.68
$ws.Range("C7").Value = 30871
$ws.Range("D7").Value = 564570.27

# --- Sheet: Regular ---
$ws = $wb.Worksheets.Item("Regular")
$ws.Range("C2").Value = 1506926
$ws.Range("D2").Value = 3864197.5
$ws.Range("C3").Value = 630653
$ws.Range("D3").Value = 1753631.14
$ws.Range("C4").Value = 1817418
$ws.Range("D4").Value = 4692244.73
$ws.Range("E4").Value = 2531662
$ws.Range("C5").Value = 762316
$ws.Range("D5").Value = 2086748.89
$ws.Range("C6").Value = 2320330
$ws.Range("D6").Value = 6323928.79
$ws.Range("E6").Value = 3227981
$ws.Range("C7").Value = 792865
$ws.Range("D7").Value = 2264608.14

# --- Sheet: Others ---
$ws = $wb.Worksheets.Item("Others")
$ws.Range("C2").Value = 12023
$ws.Range("D2").Value = 163111.66
$ws.Range("C3").Value = 5257
$ws.Range("D3").Value = 74651.24000000001
$ws.Range("C4").Value = 14373
$ws.Range("D4").Value = 189633.38
$ws.Range("C5").Value = 5425
$ws.Range("D5").Value = 75294.38
$ws.Range("C6").Value = 20288
$ws.Range("D6").Value = 285533.96
$ws.Range("C7").Value = 6868
$ws.Range("D7").Value = 101585.84

$wb.Save()
